$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9's "blue" column (G9) changes from "yes" to "confirm".
$ws.Range("G9").Value = "confirm"

# Update the sheet's active cell / selection from G7 to G10.
$null = $ws.Range("G10").Select()
